$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 552.5393490234258
$ws.Range("D2").Value = 135.9824027452931
$ws.Range("F2").Value = 450
$ws.Range("G2").Value = 509
$ws.Range("H2").Value = 621

$ws.Range("C3").Value = 40.34784979735605
$ws.Range("D3").Value = 5.33264974160819
$ws.Range("F3").Value = 37.23
$ws.Range("G3").Value = 40.07
$ws.Range("H3").Value = 43.47

$ws.Range("C4").Value = 1.63208194363162
$ws.Range("D4").Value = 2.160765135266668
$ws.Range("F4").Value = 0.5600000000000001
$ws.Range("G4").Value = 1.09
$ws.Range("H4").Value = 2.07

$ws.Range("C5").Value = 323.4778055602654
$ws.Range("D5").Value = 11.06650058287537
$ws.Range("F5").Value = 316.5
$ws.Range("G5").Value = 324.78
$ws.Range("H5").Value = 332.27

$ws.Range("C6").Value = 20.78425754650419
$ws.Range("D6").Value = 2.542140194760566
$ws.Range("F6").Value = 19.4
$ws.Range("H6").Value = 22.26

$ws.Range("C7").Value = -76.08580726674069
$ws.Range("D7").Value = 22.893064651485
$ws.Range("F7").Value = -92

$ws.Range("C8").Value = 7.642853545527799
$ws.Range("D8").Value = 6.897861818145037

$ws.Range("C9").Value = 9.321505417512157
$ws.Range("D9").Value = 1.685257834954974

$ws.Range("C10").Value = 867.8301581224396
$ws.Range("D10").Value = 0.4614251634769738

$ws.Range("C11").Value = 0.5554575355127357
$ws.Range("D11").Value = 0.5887772280295234

$ws.Range("C12").Value = 22.74046260795371
$ws.Range("D12").Value = 12.29178319313078

$ws.Range("C13").Value = 0.6738663151619237
$ws.Range("D13").Value = 0.7505432893694328

$ws.Range("C14").Value = 1.827044351949593
$ws.Range("D14").Value = 1.66412783741639

$ws.Range("C15").Value = 93.48580726674052
$ws.Range("D15").Value = 22.893064651485
$ws.Range("H15").Value = 109.4

$ws.Range("C16").Value = -85.32693461185355
$ws.Range("D16").Value = 20.60411647482207
$ws.Range("F16").Value = -101.9574620641016
$ws.Range("G16").Value = -83.69305820175224
$ws.Range("H16").Value = -67.69305820175224

$ws.Range("C17").Value = -77.68408106632575
$ws.Range("D17").Value = 25.30347243440544
$ws.Range("F17").Value = -92.71081852649533
$ws.Range("G17").Value = -72.79009749652566
$ws.Range("H17").Value = -57.21238401914255
